$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.114.03'
$ws.Range("E2").Value = '  +1.16%  '

$ws.Range("D3").Value = '1.832.46'
$ws.Range("E3").Value = '  +1.29%  '

$ws.Range("E4").Value = '  +0.77%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.43'
$ws.Range("E5").Value = '  +3.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.616'
$ws.Range("E6").Value = '  +1.84%  '

$ws.Range("E7").Value = '  +0.64%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.85'
$ws.Range("E8").Value = '  +5.85%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.310'
$ws.Range("E9").Value = '  +6.17%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0690'
$ws.Range("E10").Value = '  +2.46%  '

$ws.Range("E11").Value = '  +0.41%  '

$ws.Range("D12").Value = '2.102.33'
$ws.Range("E12").Value = '  +1.64%  '

$ws.Range("D13").Value = '1.851.39'
$ws.Range("E13").Value = '  +2.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '11.21'
$ws.Range("E14").Value = '  +2.63%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.665'
$ws.Range("E15").Value = '  +4.73%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '4.69'
$ws.Range("E16").Value = '  +6.75%  '

$ws.Range("D17").Value = '35.129.52'
$ws.Range("E17").Value = '  +1.31%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '70.27'
$ws.Range("E18").Value = '  +3.32%  '

$ws.Range("D19").Value = '0.0₃0792'
$ws.Range("E19").Value = '  +2.91%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '240.23'
$ws.Range("E20").Value = '  -0.52%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.88'
$ws.Range("E21").Value = '  +7.08%  '

$ws.Range("E22").Value = '  +0.51%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.59'
$ws.Range("E23").Value = '  +12.07%  '

$ws.Range("E24").Value = '  +4.11%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '171.43'
$ws.Range("E25").Value = '  -0.28%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.82'
$ws.Range("E26").Value = '  +1.23%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.55'
$ws.Range("E27").Value = '  +0.09%  '

$ws.Range("E28").Value = '  -0.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.59'
$ws.Range("E29").Value = '  +29.90%  '

$ws.Range("E30").Value = '  +0.71%  '

$ws.Range("D31").Value = '3.343.28'
$ws.Range("E31").Value = '  +37.60%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0555'
$ws.Range("E32").Value = '  +7.73%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.92'
$ws.Range("E33").Value = '  +3.65%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.02'
$ws.Range("E34").Value = '  +4.33%  '

$ws.Range("E35").Value = '  -0.46%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '93.33'
$ws.Range("E36").Value = '  +10.55%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.681'
$ws.Range("E37").Value = '  +6.00%  '

$ws.Range("E38").Value = '  +6.18%  '

$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '1.326.18'
$ws.Range("E39").Value = '  +1.33%  '

$ws.Range("B40").Value = 'VeChain'
$ws.Range("C40").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0194'
$ws.Range("E40").Value = '  +2.95%  '

$ws.Range("B41").Value = 'InjectiveProtocol'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '15.04'
$ws.Range("E41").Value = '  +0.61%  '

$ws.Range("E42").Value = '  +2.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.995'
$ws.Range("E43").Value = '  +5.23%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.37'
$ws.Range("E44").Value = '  +0.83%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.45'
$ws.Range("E45").Value = '  +1.05%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.79'
$ws.Range("E46").Value = '  +0.00%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.23'
$ws.Range("E47").Value = '  +8.51%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0507'
$ws.Range("E48").Value = '  -2.57%  '

$ws.Range("D49").Value = '2.009.86'
$ws.Range("E49").Value = '  +2.13%  '

$ws.Range("E50").Value = '  +0.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '100.93'
$ws.Range("E51").Value = '  -0.28%  '
